# Insert a new column "is_normal_for_donor" (with value "No") right after
# the existing "is_normal" column (currently column C), pushing the
# remaining columns (Sample_ID, test_col_2ID2, relative_file_path) one to
# the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns D:F -> E:G by inserting a new column before column D.
$ws.Columns.Item(4).Insert()

# New column D header + value.
$ws.Cells.Item(1, 4).Value = "is_normal_for_donor"
$ws.Cells.Item(2, 4).Value = "No"

# Match the target column width for the newly inserted column (matches
# neighbouring column C's width, but without the "best fit" autosize flag).
$ws.Columns.Item(4).ColumnWidth = 8.33

# Update the selection to match the target state.
$ws.Range("E6").Select()
